$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Archivos cell in row 7 to the new value "style.css"
$ws.Range("C7").Value = "style.css"

# Update the active cell selection to C8, matching the saved view state
$ws.Range("C8").Select()
